$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 849.37256
$ws.Range("I15").Value = 849.37256
$ws.Range("K15").Value = 2548.11768
$ws.Range("M15").Value = -2379.11768
$ws.Range("H48").Value = 2319
$ws.Range("J48").Value = 2319
$ws.Range("L48").Value = 6957
$ws.Range("N48").Value = -7541
$ws.Range("H56").Value = 2319
$ws.Range("J56").Value = 2319
$ws.Range("L56").Value = 6957
$ws.Range("N56").Value = -8025
$ws.Range("H62").Value = 500000500
$ws.Range("I62").Value = 500000500
$ws.Range("K62").Value = 500000500
$ws.Range("M62").Value = -499999876
$ws.Range("H65").Value = 500000500
$ws.Range("I65").Value = 500000500
$ws.Range("K65").Value = 2500002500
$ws.Range("M65").Value = -2499999380
$ws.Range("H94").Value = 2336.2856
$ws.Range("I94").Value = 2336.2856
$ws.Range("K94").Value = 2336.2856
$ws.Range("M94").Value = -1885.2856
$ws.Range("H101").Value = 1482134.4
$ws.Range("I101").Value = 2222389
$ws.Range("J101").Value = 1625
$ws.Range("K101").Value = 6667167
$ws.Range("L101").Value = 4875
$ws.Range("M101").Value = -6665545
$ws.Range("N101").Value = -8119
$ws.Range("H131").Value = 1405.7858
$ws.Range("I131").Value = 653.6667
$ws.Range("K131").Value = 1961.0001
$ws.Range("M131").Value = 3078.9999
$ws.Range("H132").Value = 1042.3125
$ws.Range("I132").Value = 925.6799999999999
$ws.Range("J132").Value = 1458.8572
$ws.Range("K132").Value = 2777.04
$ws.Range("L132").Value = 4376.571599999999
$ws.Range("M132").Value = -247.04
$ws.Range("N132").Value = -9436.571599999999
$ws.Range("H137").Value = 1755.75
$ws.Range("I137").Value = 1191
$ws.Range("K137").Value = 3573
$ws.Range("M137").Value = -1023
$ws.Range("H140").Value = 81148.64
$ws.Range("J140").Value = 81148.64
$ws.Range("L140").Value = 81148.64
$ws.Range("N140").Value = -91508.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 266.33334
$ws.Range("I5").Value = 266.33334
$ws.Range("K5").Value = 266.33334
$ws.Range("M5").Value = -154.33334
$ws.Range("H32").Value = 5548.39
$ws.Range("I32").Value = 4615.2563
$ws.Range("K32").Value = 4615.2563
$ws.Range("M32").Value = -4328.2563
$ws.Range("H61").Value = 5029.6553
$ws.Range("I61").Value = 6281.4116
$ws.Range("J61").Value = 3256.3333
$ws.Range("K61").Value = 6281.4116
$ws.Range("L61").Value = 3256.3333
$ws.Range("M61").Value = -6069.4116
$ws.Range("N61").Value = -3680.3333
$ws.Range("H74").Value = 1629.7931
$ws.Range("I74").Value = 455.5
$ws.Range("J74").Value = 4239.3335
$ws.Range("K74").Value = 455.5
$ws.Range("L74").Value = 4239.3335
$ws.Range("M74").Value = 418.5
$ws.Range("N74").Value = -5987.3335
$ws.Range("H77").Value = 1629.7931
$ws.Range("I77").Value = 455.5
$ws.Range("J77").Value = 4239.3335
$ws.Range("K77").Value = 2277.5
$ws.Range("L77").Value = 21196.6675
$ws.Range("M77").Value = 2090.5
$ws.Range("N77").Value = -29932.6675
$ws.Range("H102").Value = 1408.5555
$ws.Range("I102").Value = 1329
$ws.Range("J102").Value = 1448.3334
$ws.Range("K102").Value = 1329
$ws.Range("L102").Value = 1448.3334
$ws.Range("M102").Value = 293
$ws.Range("N102").Value = -4692.3334
$ws.Range("H132").Value = 1621.7576
$ws.Range("I132").Value = 1484.0333
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 4452.0999
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -1922.0999
$ws.Range("N132").Value = -14057
$ws.Range("H136").Value = 5029.6553
$ws.Range("I136").Value = 6281.4116
$ws.Range("J136").Value = 3256.3333
$ws.Range("K136").Value = 18844.2348
$ws.Range("L136").Value = 9768.999899999999
$ws.Range("M136").Value = -16294.2348
$ws.Range("N136").Value = -14868.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 266.33334
$ws.Range("I4").Value = 266.33334
$ws.Range("K4").Value = 266.33334
$ws.Range("M4").Value = -151.33334
$ws.Range("H86").Value = 86093.914
$ws.Range("I86").Value = 3104.7778
$ws.Range("J86").Value = 335061.34
$ws.Range("K86").Value = 3104.7778
$ws.Range("L86").Value = 335061.34
$ws.Range("M86").Value = -1981.7778
$ws.Range("N86").Value = -337307.34
$ws.Range("H89").Value = 86093.914
$ws.Range("I89").Value = 3104.7778
$ws.Range("J89").Value = 335061.34
$ws.Range("K89").Value = 15523.889
$ws.Range("L89").Value = 1675306.7
$ws.Range("M89").Value = -9907.888999999999
$ws.Range("N89").Value = -1686538.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1277.6923
$ws.Range("J22").Value = 1480.9
$ws.Range("L22").Value = 1480.9
$ws.Range("N22").Value = -2180.9
$ws.Range("H26").Value = 9921
$ws.Range("J26").Value = 9921
$ws.Range("L26").Value = 9921
$ws.Range("N26").Value = -10495
$ws.Range("H31").Value = 2651.818
$ws.Range("I31").Value = 1829.3334
$ws.Range("J31").Value = 2960.25
$ws.Range("K31").Value = 1829.3334
$ws.Range("L31").Value = 2960.25
$ws.Range("M31").Value = -1534.3334
$ws.Range("N31").Value = -3550.25
$ws.Range("H34").Value = 2651.818
$ws.Range("I34").Value = 1829.3334
$ws.Range("J34").Value = 2960.25
$ws.Range("K34").Value = 1829.3334
$ws.Range("L34").Value = 2960.25
$ws.Range("M34").Value = -1627.3334
$ws.Range("N34").Value = -3364.25
$ws.Range("H58").Value = 1861.25
$ws.Range("I58").Value = 1218.4375
$ws.Range("K58").Value = 1218.4375
$ws.Range("M58").Value = -1015.4375
$ws.Range("H136").Value = 1861.25
$ws.Range("I136").Value = 1218.4375
$ws.Range("K136").Value = 3655.3125
$ws.Range("M136").Value = -1105.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 324
$ws.Range("I18").Value = 198.66667
$ws.Range("K18").Value = 596.00001
$ws.Range("M18").Value = -427.00001
$ws.Range("H113").Value = 825.35
$ws.Range("I113").Value = 989.3333
$ws.Range("J113").Value = 796.41174
$ws.Range("K113").Value = 2967.9999
$ws.Range("L113").Value = 2389.23522
$ws.Range("M113").Value = -797.9998999999998
$ws.Range("N113").Value = -6729.23522
$ws.Range("H131").Value = 17434.785
$ws.Range("J131").Value = 19195.553
$ws.Range("L131").Value = 57586.659
$ws.Range("N131").Value = -67666.659
$ws.Range("H133").Value = 4750
$ws.Range("I133").Value = 3500
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 10500
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -5440
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2557.125
$ws.Range("I80").Value = 2951.4
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 2951.4
$ws.Range("L80").Value = 1900
$ws.Range("M80").Value = -1953.4
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 2557.125
$ws.Range("I83").Value = 2951.4
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 14757
$ws.Range("L83").Value = 9500
$ws.Range("M83").Value = -9765
$ws.Range("N83").Value = -19484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3500
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 3500
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488
$ws.Range("H82").Value = 1761
$ws.Range("I82").Value = 1521.875
$ws.Range("J82").Value = 2034.2858
$ws.Range("K82").Value = 1521.875
$ws.Range("L82").Value = 2034.2858
$ws.Range("M82").Value = -1160.875
$ws.Range("N82").Value = -2756.2858
$ws.Range("H85").Value = 1761
$ws.Range("I85").Value = 1521.875
$ws.Range("J85").Value = 2034.2858
$ws.Range("K85").Value = 1521.875
$ws.Range("L85").Value = 2034.2858
$ws.Range("M85").Value = -273.875
$ws.Range("N85").Value = -4530.2858
$ws.Range("H132").Value = 1755.5625
$ws.Range("I132").Value = 1504.4736
$ws.Range("J132").Value = 2122.5386
$ws.Range("K132").Value = 4513.4208
$ws.Range("L132").Value = 6367.6158
$ws.Range("M132").Value = -1983.4208
$ws.Range("N132").Value = -11427.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2327
$ws.Range("I132").Value = 1576.4667
$ws.Range("K132").Value = 4729.4001
$ws.Range("M132").Value = -2199.4001
$ws.Range("H136").Value = 1302.8064
$ws.Range("I136").Value = 884.2381
$ws.Range("K136").Value = 2652.7143
$ws.Range("M136").Value = -102.7143000000001
